# Update generated-data counts ("想去人数" / want-to-go counts) across sheets.
# Mirrors the upstream gh-pages data regeneration: a handful of F-column
# (column 6) numeric values bump by a small delta on the 展览 (Exhibits) and
# 演出 (Shows) sheets, plus the same rows mirrored into 全部类型 (All types).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 810
$wsExpo.Range("F4").Value = 13532
$wsExpo.Range("F5").Value = 13355
$wsExpo.Range("F6").Value = 1034
$wsExpo.Range("F13").Value = 714
$wsExpo.Range("F15").Value = 41

# 演出 (Shows) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 21
$wsShow.Range("F8").Value = 825

# 全部类型 (All types) sheet - combined listing, same rows shifted down
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 810
$wsAll.Range("F6").Value = 13532
$wsAll.Range("F7").Value = 13355
$wsAll.Range("F8").Value = 1034
$wsAll.Range("F15").Value = 714
$wsAll.Range("F17").Value = 21
$wsAll.Range("F19").Value = 41
$wsAll.Range("F33").Value = 825
